$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Kartik Tyagi, Rajasthan Royals): runs 2 -> 0, balls 3 -> 0
$ws.Range("C2").Value = "'0"
$ws.Range("D2").Value = "'0"

# Row 4 (Kartik Tyagi, Rajasthan Royals): runs 0 -> 2, balls 0 -> 3
$ws.Range("C4").Value = "'2"
$ws.Range("D4").Value = "'3"
